# antibiotic_names.xlsx — add two new antibiotic-name mapping rows
# (mirrors the "Add files via upload" commit):
#   - "nafcillin" / "nafcillin"               inserted after "moxifloxacin"
#   - "quinupristin + dalfopristin" /
#     "quinupristin dalfopristin"              inserted after "plazomicin"
# and keep the AutoFilter range / _FilterDatabase defined name in sync with
# the now-larger table (A1:B121 -> A1:B123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "nafcillin" row --------------------------------------------
# Before the edit, row 92 is "nalidixic acid" (right after "moxifloxacin"
# in row 91). Inserting a whole row here pushes everything down by one and
# leaves a blank row 92 for the new entry.
$ws.Rows.Item(92).EntireRow.Insert()
$ws.Range("A92").Value2 = "nafcillin"
$ws.Range("B92").Value2 = "nafcillin"

# --- Insert "quinupristin + dalfopristin" row --------------------------
# After the first insert, "plazomicin" has shifted down to row 103 and
# "rifampin" now sits at row 104. Insert a new row 104 for the addition.
$ws.Rows.Item(104).EntireRow.Insert()
$ws.Range("A104").Value2 = "quinupristin + dalfopristin"
$ws.Range("B104").Value2 = "quinupristin dalfopristin"

# --- Re-apply the AutoFilter over the new A1:B123 extent ---------------
# Excel's AutoFilter toggles off if re-invoked while already active, so
# switch it off first, then turn it back on over the full updated range.
$ws.AutoFilterMode = $False
$ws.Range("A1:B123").AutoFilter()

# --- Keep the _FilterDatabase defined name pointing at the full table --
foreach ($n in $wb.Names) {
    if ($n.Name -eq "antibiotic_names!_FilterDatabase") {
        $n.RefersTo = "=antibiotic_names!`$A`$1:`$B`$123"
    }
}

# --- Match the author's final selection in the sheet view --------------
$ws.Range("E103").Select()
